# Weekly update: insert a new week's Primera/Segunda price records at the
# top of the data block (rows 430-431), pushing all subsequent rows down
# by two, and extending the table with the two rows that fall off the
# bottom of the old range (now rows 548-549).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Snapshot the rows that will be pushed past the old last row (547)
#    before anything else is overwritten.
$tailBlock = $ws.Range("A546:R547").Value()

# 2) Shift the existing data block (rows 430-545) down by two rows, into
#    432-547. Captured first so source/destination overlap is safe.
$shiftBlock = $ws.Range("A430:R545").Value()
$ws.Range("A432:R547").Value = $shiftBlock

# 3) Write the new week's two quality rows (Primera / Segunda) into the
#    now-vacated rows 430-431.
$ws.Range("D430").Value = 44642
$ws.Range("I430").Value = "Primera"
$ws.Range("J430").Value = 1130
$ws.Range("K430").Value = 1500
$ws.Range("L430").Value = 1600
$ws.Range("M430").Value = 1540
$ws.Range("O430").Value = "Provincia de Quillota"
$ws.Range("P430").Value = 1540

$ws.Range("D431").Value = 44642
$ws.Range("I431").Value = "Segunda"
$ws.Range("J431").Value = 650
$ws.Range("K431").Value = 1200
$ws.Range("L431").Value = 1200
$ws.Range("M431").Value = 1200
$ws.Range("O431").Value = "Provincia de Quillota"
$ws.Range("P431").Value = 1200

# 4) Append the two rows captured in step 1 as the new rows 548-549.
$ws.Range("A548:R549").Value = $tailBlock
